$wb = $excel.ActiveWorkbook

# Helper pattern used below for the handful of cells whose text looks like a
# calendar date ("2026-01-28"): assigning that text straight to .Value makes
# Excel auto-convert it into a real date serial, which is NOT what the
# source data uses (every date/time in this workbook - e.g. "2026-01-27" on
# the Reports sheet - is stored as literal text). Building it as a formula
# that RETURNS the text, then copy/paste-special-as-values over itself,
# freezes it back down to a plain text cell without ever invoking the
# date parser, and without leaving any stray number-format behind.
function Set-TextValue {
    param($range, [string]$text)
    $escaped = $text.Replace("""", """""")
    $range.Formula = "=""" + $escaped + """"
    $range.Copy()
    $range.PasteSpecial("xlPasteValues")
}

# ---------------------------------------------------------------------------
# Sheet "Scans": append a new scan-log row (row 6)
# ---------------------------------------------------------------------------
$scans = $wb.Worksheets.Item("Scans")

Set-TextValue $scans.Range("A6") "2026-01-28"
$scans.Range("B6").Value = "00:21:31"

# Grow the "Scans" table to include the new row.
$scansTable = $scans.ListObjects.Item("Scans")
$scansTable.Resize($scans.Range("A1:C6"))

# ---------------------------------------------------------------------------
# Sheet "Order Papers": append two new rows (7 and 8) for newly-missing
# Order Paper entries.
# ---------------------------------------------------------------------------
$op = $wb.Worksheets.Item("Order Papers")

# Row 7 - Public Accounts / 63rd Report
Set-TextValue $op.Range("A7") "2026-01-28"
$op.Range("B7").Value = "Public Accounts"
$op.Range("C7").Value = "63rd Report: Increasing police productivity"
$op.Range("D7").Value = "HC 1239"
Set-TextValue $op.Range("E7") "2026-01-28"
$op.Range("F7").Value = "00:01:00"
$op.Range("G7").Value = "Missing"

# Row 8 - Transport / 5th Report
Set-TextValue $op.Range("A8") "2026-01-28"
$op.Range("B8").Value = "Transport"
$op.Range("C8").Value = "5th Report: Engine for growth: securing skills for transport manufacturing"
$op.Range("D8").Value = "HC 1223"
Set-TextValue $op.Range("E8") "2026-01-28"
$op.Range("F8").Value = "00:01:00"
$op.Range("G8").Value = "Missing"

# Grow the "Order_Papers" table to include the two new rows.
$opTable = $op.ListObjects.Item("Order_Papers")
$opTable.Resize($op.Range("A1:H8"))

$excel.CutCopyMode = 0
